# Trade #5 closed at 2026-02-16 22:57:43 - base_strategy UP +0.000%
# Append the new trade row (row 6) to both the "All Trades" and
# "base_strategy" sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(6, 1).Value = 5

    # Force text storage for the date-looking string so it is not
    # auto-converted to a date serial number, then restore the default
    # (unstyled) cell style so no stray NumberFormat sticks around.
    $ws.Cells.Item(6, 2).NumberFormat = "@"
    $ws.Cells.Item(6, 2).Value = "2026-02-16"
    $ws.Cells.Item(6, 2).Style = "Normal"

    $ws.Cells.Item(6, 3).Value = "22:57:43"
    $ws.Cells.Item(6, 4).Value = "base_strategy"
    $ws.Cells.Item(6, 5).Value = "UP"
    $ws.Cells.Item(6, 6).Value = 0.5
    $ws.Cells.Item(6, 7).Value = ""
    $ws.Cells.Item(6, 8).Value = "OPEN"
    $ws.Cells.Item(6, 9).Value = 0
    $ws.Cells.Item(6, 10).Value = 0
    $ws.Cells.Item(6, 11).Value = 100
    $ws.Cells.Item(6, 12).Value = 0
    $ws.Cells.Item(6, 13).Value = 0
    $ws.Cells.Item(6, 14).Value = 0.6
    $ws.Cells.Item(6, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(6, 16).Value = ""
    $ws.Cells.Item(6, 17).Value = 0
}
